$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Row 4: the "Online" format now also applies to columns E and F ---
$ws.Range("E4").Value = $ws.Range("D4").Value2
$ws.Range("F4").Value = $ws.Range("D4").Value2

# --- Row 6: the slides link now points at the handout PDF ---
$ws.Range("O6").Value = "EIWG_Clinical_Pharmacology_PSI_poster_02June2022_handout.pdf"

# --- Rows 22 & 25: the speaker/company pairs were swapped between the two rows ---
$l22 = $ws.Range("L22").Value2
$m22 = $ws.Range("M22").Value2
$l25 = $ws.Range("L25").Value2
$m25 = $ws.Range("M25").Value2

$ws.Range("L22").Value = $l25
$ws.Range("M22").Value = $m25
$ws.Range("L25").Value = $l22
$ws.Range("M25").Value = $m22

# --- View state: scroll the frozen pane down and select L22:M22 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$ws.Range("L22:M22").Select()
